$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 9 (pushes the old row 9 "Person FBI Identification ID"
# under Post Consolidation Identifiers down to row 10), inheriting formatting from
# the row above it.
$ws.Rows.Item(9).Insert()

$ws.Range("A9").Value = "Person State Fingerprint ID"
$ws.Range("B9").Value = "An identification of a person based on a Fingerprint ID."
$ws.Range("C9").Value = "/CHcr-doc:CriminalHistoryConsolidationReport/nc:Person/CHcr-ext:PostConsolidationIdentifiers/j:PersonStateFingerprintIdentification/nc:IdentificationID"
